$wb = $excel.ActiveWorkbook

# Both "展览" (sheet1) and "全部类型" (sheet4) contain the same data table
# and need their "想去人数" (column F) values updated.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 2152
    $ws.Range("F4").Value = 1538
    $ws.Range("F5").Value = 7252
    $ws.Range("F7").Value = 166
}
